$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.993.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.882.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4776"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3952"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08033"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.020"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.881.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.069"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.213"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.014"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06749"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.984.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.526"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.104.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.109"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.504"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9789"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09587"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.637"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.347"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.361"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06078"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02249"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.209"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.197"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.011"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6005"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1902"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.263"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5689"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.936"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.348"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06807"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.35%  "
